$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 2 (old orders get pushed down to rows 4-6)
$ws.Rows("2:3").Insert()

# New order (was Bill No 437) in row 2
$ws.Cells.Item(2, 1).Value = 437
$ws.Cells.Item(2, 2).Value = 45725.22928240741
$ws.Cells.Item(2, 2).NumberFormat = "m/d/yy"
$ws.Cells.Item(2, 3).Value = "Ajay Francis Anchan"
$ws.Cells.Item(2, 4).Value = 9
$ws.Cells.Item(2, 5).Value = 100
$ws.Cells.Item(2, 6).Value = 1.3
$ws.Cells.Item(2, 7).Value = 1.3
$ws.Cells.Item(2, 8).Value = 2.6
$ws.Cells.Item(2, 9).Value = "Butterscotch Lassi (x1), 8PM Coffee (x1)"

# New order (was Bill No 436) in row 3
$ws.Cells.Item(3, 1).Value = 436
$ws.Cells.Item(3, 2).Value = 45725.22928240741
$ws.Cells.Item(3, 2).NumberFormat = "m/d/yy"
$ws.Cells.Item(3, 3).Value = "Ajay Francis Anchan"
$ws.Cells.Item(3, 4).Value = 8
$ws.Cells.Item(3, 5).Value = 160
$ws.Cells.Item(3, 6).Value = 2.2
$ws.Cells.Item(3, 7).Value = 2.2
$ws.Cells.Item(3, 8).Value = 4.4
$ws.Cells.Item(3, 9).Value = "Butterscotch Lassi (x2), 8PM Coffee (x1)"

Write-Host "done"
